$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.603.77'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.881.92'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.88%  '
$ws.Range("D7").Value = '3.878.30'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.85%  '
$ws.Range("D15").Value = '4.536.58'
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '3.888.43'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '69.606.54'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.08%  '
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '490.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("E23").Value = '  +4.47%  '
$ws.Range("E24").Value = '  +3.72%  '
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("D32").Value = '4.032.13'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").Value = '3.844.69'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.57%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.141'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +3.20%  '
$ws.Range("E43").Value = '  +6.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '436.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.99%  '
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000274'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +21.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("E51").Value = '  +4.13%  '
